# Update the "Generate Report for Handback" timestamps.
#
# Overview!G2 ("Latest HO Xliff Generate Date") and de-de!H2
# ("Correspond Handoff Datetime") both originally held the same
# timestamp text, so both move to the new handoff time.
#
# zh-cn!H2/K2 and de-de!K2 hold the per-locale handoff/handback
# timestamps that also need to be bumped forward.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview) / Correspond Handoff Datetime (de-de)
$wsOverview.Range("G2").Value = "2016-08-19 17:08:28"
$wsDeDe.Range("H2").Value = "2016-08-19 17:08:28"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-19 17:08:23"
$wsZhCn.Range("K2").Value = "2016-08-19 17:08:40"

# de-de: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-19 17:08:47"
